$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gphb5"
$ws.Cells.Item(2, 3).Value = "Tshr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.3862496666666667
$ws.Cells.Item(2, 8).Value = 1.158749
$ws.Cells.Item(2, 9).Value = 0.06565390438652881
$ws.Cells.Item(2, 10).Value = 0.06565390438652881
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.960636
$ws.Cells.Item(2, 14).Value = 2.881908
$ws.Cells.Item(2, 15).Value = 0.2124324572954377
$ws.Cells.Item(2, 16).Value = 0.2124324572954377
$ws.Cells.Item(2, 17).Value = 0.371045334788
$ws.Cells.Item(2, 18).Value = 3.339408013092
$ws.Cells.Item(2, 19).Value = 0.01394702023987003
$ws.Cells.Item(2, 20).Value = 0.01394702023987003

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gphb5"
$ws.Cells.Item(3, 3).Value = "Tshr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.3862496666666667
$ws.Cells.Item(3, 8).Value = 1.158749
$ws.Cells.Item(3, 9).Value = 0.06565390438652881
$ws.Cells.Item(3, 10).Value = 0.06565390438652881
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.5285266666666667
$ws.Cells.Item(3, 14).Value = 1.58558
$ws.Cells.Item(3, 15).Value = 0.116876963330717
$ws.Cells.Item(3, 16).Value = 0.116876963330717
$ws.Cells.Item(3, 17).Value = 0.2041432488244445
$ws.Cells.Item(3, 18).Value = 1.83728923942
$ws.Cells.Item(3, 19).Value = 0.00767342897550273
$ws.Cells.Item(3, 20).Value = 0.00767342897550273

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gphb5"
$ws.Cells.Item(4, 3).Value = "Tshr"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.3862496666666667
$ws.Cells.Item(4, 8).Value = 1.158749
$ws.Cells.Item(4, 9).Value = 0.06565390438652881
$ws.Cells.Item(4, 10).Value = 0.06565390438652881
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.032914666666667
$ws.Cells.Item(4, 14).Value = 9.098744
$ws.Cells.Item(4, 15).Value = 0.6706905793738452
$ws.Cells.Item(4, 16).Value = 0.6706905793738454
$ws.Cells.Item(4, 17).Value = 1.171462279028444
$ws.Cells.Item(4, 18).Value = 10.543160511256
$ws.Cells.Item(4, 19).Value = 0.04403345517115605
$ws.Cells.Item(4, 20).Value = 0.04403345517115605

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gphb5"
$ws.Cells.Item(5, 3).Value = "Tshr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.711678
$ws.Cells.Item(5, 8).Value = 14.135034
$ws.Cells.Item(5, 9).Value = 0.8008810974044716
$ws.Cells.Item(5, 10).Value = 0.8008810974044716
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.960636
$ws.Cells.Item(5, 14).Value = 2.881908
$ws.Cells.Item(5, 15).Value = 0.2124324572954377
$ws.Cells.Item(5, 16).Value = 0.2124324572954377
$ws.Cells.Item(5, 17).Value = 4.526207507208
$ws.Cells.Item(5, 18).Value = 40.735867564872
$ws.Cells.Item(5, 19).Value = 0.1701331395230987
$ws.Cells.Item(5, 20).Value = 0.1701331395230987

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gphb5"
$ws.Cells.Item(6, 3).Value = "Tshr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4.711678
$ws.Cells.Item(6, 8).Value = 14.135034
$ws.Cells.Item(6, 9).Value = 0.8008810974044716
$ws.Cells.Item(6, 10).Value = 0.8008810974044716
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.5285266666666667
$ws.Cells.Item(6, 14).Value = 1.58558
$ws.Cells.Item(6, 15).Value = 0.116876963330717
$ws.Cells.Item(6, 16).Value = 0.116876963330717
$ws.Cells.Item(6, 17).Value = 2.490247467746667
$ws.Cells.Item(6, 18).Value = 22.41222720972
$ws.Cells.Item(6, 19).Value = 0.09360455065360684
$ws.Cells.Item(6, 20).Value = 0.09360455065360684

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gphb5"
$ws.Cells.Item(7, 3).Value = "Tshr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.711678
$ws.Cells.Item(7, 8).Value = 14.135034
$ws.Cells.Item(7, 9).Value = 0.8008810974044716
$ws.Cells.Item(7, 10).Value = 0.8008810974044716
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.032914666666667
$ws.Cells.Item(7, 14).Value = 9.098744
$ws.Cells.Item(7, 15).Value = 0.6706905793738452
$ws.Cells.Item(7, 16).Value = 0.6706905793738454
$ws.Cells.Item(7, 17).Value = 14.29011731081067
$ws.Cells.Item(7, 18).Value = 128.611055797296
$ws.Cells.Item(7, 19).Value = 0.537143407227766
$ws.Cells.Item(7, 20).Value = 0.5371434072277661

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Gphb5"
$ws.Cells.Item(8, 3).Value = "Tshr"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.7851903333333333
$ws.Cells.Item(8, 8).Value = 2.355571
$ws.Cells.Item(8, 9).Value = 0.1334649982089996
$ws.Cells.Item(8, 10).Value = 0.1334649982089996
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.960636
$ws.Cells.Item(8, 14).Value = 2.881908
$ws.Cells.Item(8, 15).Value = 0.2124324572954377
$ws.Cells.Item(8, 16).Value = 0.2124324572954377
$ws.Cells.Item(8, 17).Value = 0.754282101052
$ws.Cells.Item(8, 18).Value = 6.788538909468
$ws.Cells.Item(8, 19).Value = 0.02835229753246897
$ws.Cells.Item(8, 20).Value = 0.02835229753246897

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Gphb5"
$ws.Cells.Item(9, 3).Value = "Tshr"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.7851903333333333
$ws.Cells.Item(9, 8).Value = 2.355571
$ws.Cells.Item(9, 9).Value = 0.1334649982089996
$ws.Cells.Item(9, 10).Value = 0.1334649982089996
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.5285266666666667
$ws.Cells.Item(9, 14).Value = 1.58558
$ws.Cells.Item(9, 15).Value = 0.116876963330717
$ws.Cells.Item(9, 16).Value = 0.116876963330717
$ws.Cells.Item(9, 17).Value = 0.4149940295755556
$ws.Cells.Item(9, 18).Value = 3.73494626618
$ws.Cells.Item(9, 19).Value = 0.01559898370160746
$ws.Cells.Item(9, 20).Value = 0.01559898370160746

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Gphb5"
$ws.Cells.Item(10, 3).Value = "Tshr"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.7851903333333333
$ws.Cells.Item(10, 8).Value = 2.355571
$ws.Cells.Item(10, 9).Value = 0.1334649982089996
$ws.Cells.Item(10, 10).Value = 0.1334649982089996
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.032914666666667
$ws.Cells.Item(10, 14).Value = 9.098744
$ws.Cells.Item(10, 15).Value = 0.6706905793738452
$ws.Cells.Item(10, 16).Value = 0.6706905793738454
$ws.Cells.Item(10, 17).Value = 2.381415278091556
$ws.Cells.Item(10, 18).Value = 21.432737502824
$ws.Cells.Item(10, 19).Value = 0.08951371697492316
$ws.Cells.Item(10, 20).Value = 0.08951371697492315
